$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Two new report rows (17 and 18 in the source report => rows 16/17 on the sheet,
# just below the existing 14 data rows), quarter 2023Q3 and 2023Q4 entries.

# --- Column A (CVR) -- new companies, stored as text like the rest of column A ---
$ws.Cells.Item(16, 1).NumberFormat = "@"
$ws.Cells.Item(16, 1).Value = "80493215"
$ws.Cells.Item(16, 1).ClearFormats()

$ws.Cells.Item(17, 1).NumberFormat = "@"
$ws.Cells.Item(17, 1).Value = "12070942"
$ws.Cells.Item(17, 1).ClearFormats()

# --- Column D (Løsning) -- reuses existing values ---
$ws.Cells.Item(16, 4).Value = "Visma Løn og HR"
$ws.Cells.Item(17, 4).Value = "Visma Løn"

# --- Column G (Ny leverandør) -- brand new suppliers ---
$ws.Cells.Item(16, 7).Value = "Lessor"
$ws.Cells.Item(17, 7).Value = "Azets"

# --- Column H (Quarter) ---
$ws.Cells.Item(16, 8).Value = "2023Q3"
$ws.Cells.Item(17, 8).Value = "2023Q4"

# --- Column I (TCV_range) -- reused for every row ---
$ws.Cells.Item(16, 9).Value = "80000-100000"
$ws.Cells.Item(17, 9).Value = "80000-100000"

# --- Column B (Year), C (TCV amount), E (Opsagt dato) ---
$ws.Cells.Item(16, 2).Value = 2023
$ws.Cells.Item(16, 3).Value = 83103.05
$ws.Cells.Item(16, 5).Value = 45198
$ws.Cells.Item(16, 5).NumberFormat = "YYYY-MM-DD HH:MM:SS"

$ws.Cells.Item(17, 2).Value = 2023
$ws.Cells.Item(17, 3).Value = 85832.52
$ws.Cells.Item(17, 5).Value = 45205
$ws.Cells.Item(17, 5).NumberFormat = "YYYY-MM-DD HH:MM:SS"
